$wb = $excel.ActiveWorkbook

# --- amd_categoryArticle: insert the new "logo" column before the meta_* columns ---
$ws = $wb.Worksheets.Item("amd_categoryArticle")

# Insert a new column at F (pushes old F..L to G..M)
$ws.Columns.Item(6).Insert()

# Fill in the new "logo" column values
$ws.Cells.Item(1, 6).Value = "logo"
$ws.Cells.Item(2, 6).Value = "nullable"
$ws.Cells.Item(3, 6).Value = "text"

# Set explicit (non bestFit) width for the new column F ~ 8.42578125 stored units
$ws.Columns.Item(6).ColumnWidth = 7.592447916666667

# Make amd_categoryArticle the active/selected sheet with its own selection
$ws.Select()
$ws.Range("I4").Select()
